# Update TPM-derived values for the Sost-Lrp4 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.6795853333333334
$ws.Range("N2").Value = 2.038756
$ws.Range("O2").Value = 0.08284139605799233
$ws.Range("P2").Value = 0.08284139605799234
$ws.Range("Q2").Value = 0.034687394584
$ws.Range("R2").Value = 0.312186551256
$ws.Range("S2").Value = 0.08284139605799233
$ws.Range("T2").Value = 0.08284139605799234

# Row 3
$ws.Range("O3").Value = 0.3136748993401273
$ws.Range("P3").Value = 0.3136748993401273
$ws.Range("S3").Value = 0.3136748993401273
$ws.Range("T3").Value = 0.3136748993401273

# Row 4
$ws.Range("M4").Value = 2.840162333333334
$ws.Range("N4").Value = 8.520487000000001
$ws.Range("O4").Value = 0.3462155540800247
$ws.Range("P4").Value = 0.3462155540800247
$ws.Range("Q4").Value = 0.144967565818
$ws.Range("R4").Value = 1.304708092362
$ws.Range("S4").Value = 0.3462155540800247
$ws.Range("T4").Value = 0.3462155540800247

# Row 5
$ws.Range("M5").Value = 2.110486666666667
$ws.Range("N5").Value = 6.33146
$ws.Range("O5").Value = 0.2572681505218555
$ws.Range("P5").Value = 0.2572681505218555
$ws.Range("Q5").Value = 0.10772346044
$ws.Range("R5").Value = 0.96951114396
$ws.Range("S5").Value = 0.2572681505218555
$ws.Range("T5").Value = 0.2572681505218555
